# Applies the minor text corrections captured by the target diff:
#   1) Slide 4 notes ("Team Work - Skills"): "specifig" -> "specific"
#   2) Slide 6 notes ("Flexibility"): " so far or all " -> " so far or "
#
# NOTE: the diff also shows a refreshed `datetimeFigureOut` field cache
# (4/16/2018 -> 4/17/2018) in the notes master and a re-ordering of
# timestamped co-authoring entries in ppt/changesInfos/changesInfo1.xml.
# Both are PowerPoint-internal bookkeeping that is recomputed from the
# live system clock / save pipeline when a real author edits the file
# (not user-editable document content), so there is no PowerPoint
# object-model call that can reproduce them here; this script focuses on
# the actual content edit described by the commit.

$p = $ppt.ActivePresentation

# --- Slide 4 notes: "specifig" -> "specific" -------------------------------
$slide4 = $p.Slides.Item(4)
$notes4 = $slide4.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes4.Text = "Different backgrounds (Electrical Engineering, Industrial Engineering), experiences (worked for some years) and cultures, that will surely lead to controverse discussions, which will (hopefully) result in success in our project`nEverybody will state his/ her opinion and we will support each other, when somebody has minor knowledge about a specific topic`n"

# --- Slide 6 notes: " so far or all " -> " so far or " ----------------------
$slide6 = $p.Slides.Item(6)
$notes6 = $slide6.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes6.Text = "Two team members work on a topic together if possible and useful -> if one cannot contribute for any reason, we can still continue`nSome team members have experience in SCRUM, weekly meetings and short communication paths help to adapt fast to any change. Being agile can also include that we will work on topics to which we have not been used so far or even all 4 of us have to work on one taks, if neccessary. FOCUS on delivering sth after every sprint`n"
